$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1. Tidy up the selection/view on the existing "ageLE" sheet
#    (was scrolled to A4 with G17 selected; now shows the full
#    A1:D20 range selected and scrolled back to the top).
# ------------------------------------------------------------------
$ageLE = $wb.Worksheets.Item("ageLE")
[void]$ageLE.Range("A1:D20").Select()

# ------------------------------------------------------------------
# 2. Add the new "ageLE_mssa" sheet right after "ageLE" and populate
#    it with the 2021 HCAI life-expectancy age bands (3 columns:
#    lAge, uAge, ageName).
# ------------------------------------------------------------------
$newSheet = $wb.Worksheets.Add($null, $ageLE)
$newSheet.Name = "ageLE_mssa"

$newSheet.Cells.Item(1,1).Value = "lAge"
$newSheet.Cells.Item(1,2).Value = "uAge"
$newSheet.Cells.Item(1,3).Value = "ageName"
$newSheet.Cells.Item(1,3).NumberFormat = "@"

$rows = @(
    @(0,   4,   "0 - 4"),
    @(5,   9,   "5 - 9"),
    @(10,  14,  "10 - 14"),
    @(15,  19,  "15 - 19"),
    @(20,  24,  "20 - 24"),
    @(25,  29,  "25 - 29"),
    @(30,  34,  "30 - 34"),
    @(35,  39,  "35 - 39"),
    @(40,  44,  "40 - 44"),
    @(45,  49,  "45 - 49"),
    @(50,  54,  "50 - 54"),
    @(55,  59,  "55 - 59"),
    @(60,  64,  "60 - 64"),
    @(65,  69,  "65 - 69"),
    @(70,  74,  "70 - 74"),
    @(75,  79,  "75 - 79"),
    @(80,  84,  "80 - 84"),
    @(85,  199, "85 - 199")
)

$r = 2
foreach ($row in $rows) {
    $newSheet.Cells.Item($r,1).Value = $row[0]
    $newSheet.Cells.Item($r,2).Value = $row[1]
    $newSheet.Cells.Item($r,3).Value = $row[2]
    $newSheet.Cells.Item($r,3).NumberFormat = "@"
    $r = $r + 1
}

[void]$newSheet.Range("E4").Select()

# ------------------------------------------------------------------
# 3. Make the freshly-added "ageLE_mssa" sheet the active tab
#    (moves tabSelected off whichever sheet previously had it).
# ------------------------------------------------------------------
$newSheet.Activate()
